$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new literal text value (Price/Volume columns are stored
# as plain text in the workbook, so we force a Text number format while
# writing the value and then restore the default "Normal" style so the
# cell keeps its original (unstyled) appearance).
$updates = @{
    'D2' = '272.28'
    'E2' = '0.80%'
    'D3' = '26.84'
    'E3' = '0.46%'
    'D4' = '4.898'
    'E4' = '3.83%'
    'D5' = '0.06315'
    'E5' = '3.19%'
    'D6' = '6.895'
    'D7' = '3.351'
    'E7' = '5.36%'
    'D8' = '1.359'
    'E8' = '52.50%'
    'D9' = '0.8827'
    'E9' = '3.38%'
    'D10' = '0.1463'
    'E10' = '2.37%'
    'D11' = '0.05088'
    'E11' = '0.69%'
    'D12' = '0.07396'
    'E12' = '3.90%'
    'D13' = '0.03159'
    'E13' = '-0.15%'
    'D14' = '0.09029'
    'E14' = '-0.09%'
    'D15' = '0.001564'
    'E15' = '1.93%'
    'D16' = '0.0006317'
    'E16' = '3.83%'
    'D17' = '0.006021'
    'E17' = '-1.39%'
    'D18' = '3.470'
    'E18' = '0.16%'
    'E19' = '0.88%'
    'D20' = '0.3166'
    'E20' = '2.46%'
    'D21' = '0.1333'
    'E21' = '4.00%'
    'D22' = '3.900'
    'E22' = '1.48%'
    'D23' = '0.04337'
    'E23' = '2.57%'
    'D24' = '0.001179'
    'E24' = '0.19%'
    'D25' = '0.003647'
    'E25' = '-12.07%'
    'D26' = '0.0001201'
    'E26' = '-0.01%'
    'D27' = '0.0001697'
    'E27' = '0.93%'
    'E40' = '2.10%'
    'D41' = '0.006622'
    'E41' = '57.81%'
    'D42' = '0.1163'
    'E42' = '4.10%'
    'D43' = '0.002131'
    'E43' = '4.57%'
    'D44' = '0.01260'
    'E44' = '7.20%'
    'D45' = '0.00005338'
    'E45' = '4.04%'
    'E46' = '156.75%'
    'D47' = '0.02120'
    'E47' = '-29.16%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
